# account_bank_statement_import_adyen / adyen_test_credit_fees.xlsx
# Migration touch-up: the Adyen "Gross Currency" / "Net Currency" columns in
# this fixture were switched from EUR to USD, and the sheet's saved view
# (window tab split + current selection) was nudged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data change: EUR -> USD on the whole sheet (Gross Currency / Net
# Currency columns only contain "EUR" as a whole-cell value; GBP rows are
# untouched because they don't match). ---
$used = $ws.UsedRange
$replaced = $used.Replace("EUR", "USD", 1, 1, $false)

# --- View changes ---
# Shrink the sheet-tabs/horizontal-scrollbar split (tabRatio 615 -> 500).
$excel.ActiveWindow.TabRatio = 0.5

# Keep gridlines visible (matches the workbook default) while moving the
# active selection to P38, and scroll so row 10 is at the top, per the
# saved view in the updated file.
$excel.ActiveWindow.DisplayGridlines = $true
$selected = $ws.Range("P38").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
